$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 ("Give all emails a common subject prefix") is no longer needed on
# its own - delete the whole row, shifting rows 8-15 up by one.
$ws.Rows.Item(7).Delete()

# Insert a new row above the current row 15 ("Homepage") for the new
# "Cleanup" task, then fill in its value.
$ws.Rows.Item(14).Insert()
$ws.Range("B14").Value = "Remove content that is no longer used"

# Restore the selection as recorded after the edit.
$ws.Range("A7:D16").Select()
